# Update "想去人数" (column F) counts on each sheet to the refreshed
# snapshot values from the latest crawl (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7724
$ws.Range("F3").Value = 101
$ws.Range("F4").Value = 79
$ws.Range("F5").Value = 7787
$ws.Range("F8").Value = 630
$ws.Range("F11").Value = 438
$ws.Range("F12").Value = 774
$ws.Range("F14").Value = 73
$ws.Range("F15").Value = 305
$ws.Range("F16").Value = 17
$ws.Range("F17").Value = 262
$ws.Range("F18").Value = 137
$ws.Range("F19").Value = 393
$ws.Range("F23").Value = 609
$ws.Range("F24").Value = 2198
$ws.Range("F25").Value = 729
$ws.Range("F26").Value = 51
$ws.Range("F27").Value = 52
$ws.Range("F29").Value = 609
$ws.Range("F30").Value = 53

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 322
$ws.Range("F10").Value = 2

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 446

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 446
$ws.Range("F3").Value = 7724
$ws.Range("F4").Value = 101
$ws.Range("F5").Value = 79
$ws.Range("F7").Value = 7788
$ws.Range("F10").Value = 630
$ws.Range("F14").Value = 438
$ws.Range("F15").Value = 322
$ws.Range("F18").Value = 774
$ws.Range("F20").Value = 73
$ws.Range("F21").Value = 305
$ws.Range("F23").Value = 17
$ws.Range("F26").Value = 262
$ws.Range("F27").Value = 137
$ws.Range("F28").Value = 393
$ws.Range("F32").Value = 609
$ws.Range("F33").Value = 2198
$ws.Range("F34").Value = 729
$ws.Range("F35").Value = 51
$ws.Range("F36").Value = 52
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = 609
$ws.Range("F40").Value = 53
